$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("AF2").Value = 51
$ws.Range("AR2").Value = 2.5
$ws.Range("BB2").Value = 201
